$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 126, shifting existing rows 126:165 down to 127:166
$ws.Rows.Item(126).Insert()

# Populate the new row 126 with the weekly data point added by this commit
$ws.Cells.Item(126, 1).Value = 4
$ws.Cells.Item(126, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(126, 3).Value = "Los Lagos"
$ws.Cells.Item(126, 4).Value = 44463
$ws.Cells.Item(126, 5).Value = 10
$ws.Cells.Item(126, 6).Value = 100112040
$ws.Cells.Item(126, 7).Value = "Cilantro"
$ws.Cells.Item(126, 8).Value = "Sin especificar"
$ws.Cells.Item(126, 9).Value = "Primera"
$ws.Cells.Item(126, 10).Value = 300
$ws.Cells.Item(126, 11).Value = 13000
$ws.Cells.Item(126, 12).Value = 13000
$ws.Cells.Item(126, 13).Value = 13000
$ws.Cells.Item(126, 14).Value = "$/caja 36 atados"
$ws.Cells.Item(126, 15).Value = "Región Metropolitana"
$ws.Cells.Item(126, 16).Value = 361
$ws.Cells.Item(126, 17).Value = 36
$ws.Cells.Item(126, 18).Value = "Hortaliza"
